# feat: add 2022-Q1 data
#
# - Inserts a new "2022-Q1" worksheet (fund-level holdings) right before
#   the "总计" (totals) summary sheet, copying the column layout/format
#   used by the other quarterly sheets (e.g. "2021-Q4").
# - Prepends a new "2022-Q1" row to the "总计" summary table and shifts
#   the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet, positioned immediately before "总计".
#    NOTE: `Worksheets.Item(...)` is a *positional* lookup, so once the
#    sheet collection is mutated (Add/Insert/etc.) any previously
#    captured reference must be re-resolved by name rather than reused.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$formatSrc = $wb.Worksheets.Item("2021-Q4")

$q1 = $wb.Worksheets.Add($totalSheet, $null)
$q1.Name = "2022-Q1"

# Pull over the header (B1:H1) and index-column (A2:A7) look & feel from
# the "2021-Q4" sheet so the new tab matches its siblings. Re-resolve
# both sheets by name since the worksheet collection just changed.
$formatSrc = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Item("2022-Q1")

$formatSrc.Range("B1:H1").Copy($q1.Range("B1:H1"))
$formatSrc.Range("A2:A5").Copy($q1.Range("A2:A7"))

# Header labels
$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

# Fund holding rows: code, name, scale, stock position, position ratio,
# held market value (亿元), position rank. Columns B/D/E/F/G keep their
# numeric-looking values as text (matching the source data), while A
# (row index) and H (rank) are real numbers.
$q1rows = @(
    @("516970", "广发中证基建工程交易型开放式指数证券投资基金", "59.01", "99.38", "4.61", "2.7204", 8),
    @("165525", "信诚中证基建工程指数（LOF）", "17.06", "94.00", "4.37", "0.7455", 8),
    @("180020", "银华成长先锋混合", "3.05", "79.81", "4.71", "0.1437", 10),
    @("620001", "金元顺安宝石动力混合", "1.90", "56.14", "3.60", "0.0684", 5),
    @("009753", "中欧美益稳健两年持有期混合A", "2.63", "23.02", "1.43", "0.0376", 2),
    @("009754", "中欧美益稳健两年持有期混合C", "0.23", "23.02", "1.43", "0.0033", 2)
)

$r = 2
foreach ($row in $q1rows) {
    $q1.Cells.Item($r, 1).Value = $r - 2
    $q1.Cells.Item($r, 2).Value = "'" + $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = "'" + $row[2]
    $q1.Cells.Item($r, 5).Value = "'" + $row[3]
    $q1.Cells.Item($r, 6).Value = "'" + $row[4]
    $q1.Cells.Item($r, 7).Value = "'" + $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# ---------------------------------------------------------------------
# 2) Prepend the "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing quarters down by one row. Re-resolve "总计" by name again
#    since the worksheet collection changed above.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalRows = @(
    @("2022-Q1", 6, 3.72),
    @("2021-Q4", 4, 0.71),
    @("2021-Q3", 5, 5.55),
    @("2021-Q2", 7, 2.02),
    @("2021-Q1", 11, 4.01),
    @("2020-Q4", 4, 2.34)
)

$r = 2
foreach ($row in $totalRows) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]
    $r++
}
